$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.033.73'
$ws.Range("E2").Value = '  +2.06%  '

# Row 3
$ws.Range("D3").Value = '2.335.08'
$ws.Range("E3").Value = '  -0.50%  '

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = "'541.87"
$ws.Range("E5").Value = '  +5.04%  '

# Row 6
$ws.Range("D6").Value = "'135.53"
$ws.Range("E6").Value = '  +0.98%  '

# Row 8
$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = '  +0.17%  '

# Row 9
$ws.Range("D9").Value = '2.370.38'
$ws.Range("E9").Value = '  +0.69%  '

# Row 10
$ws.Range("E10").Value = '  +0.49%  '

# Row 11
$ws.Range("E11").Value = '  +0.71%  '

# Row 12
$ws.Range("D12").Value = "'5.41"
$ws.Range("E12").Value = '  +0.32%  '

# Row 13
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = '  +3.48%  '

# Row 14
$ws.Range("D14").Value = "'23.76"
$ws.Range("E14").Value = '  -1.15%  '

# Row 15
$ws.Range("D15").Value = '2.757.18'
$ws.Range("E15").Value = '  -0.18%  '

# Row 16
$ws.Range("D16").Value = '57.848.75'
$ws.Range("E16").Value = '  +1.77%  '

# Row 17
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$ws.Range("D18").Value = '2.351.40'
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").Value = "'337.71"
$ws.Range("E19").Value = '  +3.15%  '

# Row 20
$ws.Range("D20").Value = "'10.53"
$ws.Range("E20").Value = '  +0.58%  '

# Row 21
$ws.Range("E21").Value = '  +0.61%  '

# Row 22
$ws.Range("D22").Value = "'6.82"
$ws.Range("E22").Value = '  +1.40%  '

# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.10%  '

# Row 24
$ws.Range("D24").Value = "'62.02"
$ws.Range("E24").Value = '  +1.45%  '

# Row 25
$ws.Range("E25").Value = '  +1.50%  '

# Row 26
$ws.Range("D26").Value = "'8.49"
$ws.Range("E26").Value = '  -2.08%  '

# Row 27
$ws.Range("E27").Value = '  -0.45%  '

# Row 28
$ws.Range("D28").Value = "'1.38"
$ws.Range("E28").Value = '  +4.80%  '

# Row 29
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = "'174.23"
$ws.Range("E29").Value = '  +3.48%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.76"
$ws.Range("E30").Value = '  +4.24%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0740'
$ws.Range("E31").Value = '  +1.08%  '

# Row 32
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = '  -0.56%  '

# Row 33
$ws.Range("D33").Value = "'18.51"
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = '  -0.08%  '

# Row 35
$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = '  +11.39%  '

# Row 36
$ws.Range("D36").Value = "'0.996"
$ws.Range("E36").Value = '  +0.08%  '

# Row 37
$ws.Range("D37").Value = "'1.27"
$ws.Range("E37").Value = '  -1.09%  '

# Row 38
$ws.Range("D38").Value = "'4.11"
$ws.Range("E38").Value = '  +2.91%  '

# Row 39
$ws.Range("D39").Value = "'1.60"
$ws.Range("E39").Value = '  +1.88%  '

# Row 40
$ws.Range("D40").Value = "'39.30"
$ws.Range("E40").Value = '  +1.66%  '

# Row 41
$ws.Range("D41").Value = "'150.17"
$ws.Range("E41").Value = '  -0.13%  '

# Row 42
$ws.Range("E42").Value = '  -0.25%  '

# Row 43
$ws.Range("E43").Value = '  +0.80%  '

# Row 44
$ws.Range("D44").Value = "'287.43"
$ws.Range("E44").Value = '  +1.37%  '

# Row 45
$ws.Range("D45").Value = "'0.0929"
$ws.Range("E45").Value = '  -0.18%  '

# Row 46
$ws.Range("E46").Value = '  +0.65%  '

# Row 47
$ws.Range("D47").Value = "'0.564"
$ws.Range("E47").Value = '  +0.92%  '

# Row 48
$ws.Range("D48").Value = "'18.85"
$ws.Range("E48").Value = '  +2.53%  '

# Row 49
$ws.Range("E49").Value = '  +0.61%  '

# Row 50
$ws.Range("D50").Value = "'17.62"
$ws.Range("E50").Value = '  +2.11%  '

# Row 51
$ws.Range("E51").Value = '  +9.93%  '
